# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates the metrics for the "with_priors" row (row 3) in the metrics_sim sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.5540540540540541
$ws.Range("D3").Value = 0.8918918918918919

$ws.Range("H3").Value = 0.5969230769230769
$ws.Range("I3").Value = 0.09778924318324984
$ws.Range("J3").Value = 0.4594594594594595
$ws.Range("K3").Value = 125.6081081081081

$ws.Range("Q3").Value = 12
$ws.Range("R3").Value = 17
$ws.Range("S3").Value = 40
$ws.Range("T3").Value = 105
$ws.Range("U3").Value = 165
$ws.Range("V3").Value = 889
$ws.Range("W3").Value = 884
$ws.Range("X3").Value = 861
$ws.Range("Y3").Value = 796
$ws.Range("Z3").Value = 736

$ws.Range("AF3").Value = 0.986681
$ws.Range("AG3").Value = 0.981132
$ws.Range("AH3").Value = 0.955605
$ws.Range("AI3").Value = 0.883463
$ws.Range("AJ3").Value = 0.81687
